# Update "想去人数" (attendance) figures for two events that now show
# higher counts, on both the "展览" sheet and the "全部类型" sheet.
#
# 展览 (sheet1): rows 4-8 -> F4,F5,F6,F7,F8
# 全部类型 (sheet4): rows 4-6 then 9-10 -> F4,F5,F6,F9,F10

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 176
$ws1.Range("F5").Value = 3306
$ws1.Range("F6").Value = 336
$ws1.Range("F7").Value = 14
$ws1.Range("F8").Value = 418

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 176
$ws4.Range("F5").Value = 3306
$ws4.Range("F6").Value = 336
$ws4.Range("F9").Value = 14
$ws4.Range("F10").Value = 418
